$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Tipo" column from D to E, and insert a new "MAE" column at D.
$ws.Range("E1").Value = "Tipo"
$ws.Range("D1").Value = "MAE"

# Apply the header formatting (bold, centered, bordered) used by the other
# header cells to the new D1/E1 header cells.
$xlPasteFormats = -4122
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial($xlPasteFormats)
$ws.Range("C1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Set the "multiple" label for each row in the (new) Tipo column E.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 5).Value = "multiple"
}

# Updated MSE (B), R2 (C) and new MAE (D) values for each disease row.
$data = @(
    @{Row=2;  B=1.308313292598285; C=0.7868622600346196; D=0.952869278467512},
    @{Row=3;  B=6.91556435086562;  C=0.9012056313998962; D=1.96195442228898},
    @{Row=4;  B=4.071032756655834; C=0.7989074779402296; D=1.597451416511271},
    @{Row=5;  B=3.955509390133789; C=0.9974018577455562; D=1.690025121283359},
    @{Row=6;  B=2.879336475278879; C=0.9686246810746295; D=1.441611837529107},
    @{Row=7;  B=2.37101831261121; C=0.9986706441246088; D=1.265570585794963},
    @{Row=8;  B=2.010940401604862; C=0.9976965617262776; D=1.130324592094832},
    @{Row=9;  B=15.55141433873043; C=0.8140935877554698; D=3.127129605711212},
    @{Row=10; B=2.164924826073224; C=0.9935464261485367; D=1.153911192551195}
)

foreach ($row in $data) {
    $ws.Cells.Item($row.Row, 2).Value = $row.B
    $ws.Cells.Item($row.Row, 3).Value = $row.C
    $ws.Cells.Item($row.Row, 4).Value = $row.D
}

$wb.Save()
